$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.177.34"
$ws.Range("E2").Value = "  +6.75%  "
$ws.Range("D3").Value = "3.016.75"
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.39"
$ws.Range("E5").Value = "  +3.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.08"
$ws.Range("E6").Value = "  +13.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.012.01"
$ws.Range("E8").Value = "  +3.68%  "
$ws.Range("E9").Value = "  +3.83%  "
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("E11").Value = "  +5.61%  "
$ws.Range("E12").Value = "  +6.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.82"
$ws.Range("E14").Value = "  +6.82%  "
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "66.081.02"
$ws.Range("E16").Value = "  +6.65%  "
$ws.Range("D17").Value = "3.517.47"
$ws.Range("E17").Value = "  +3.74%  "
$ws.Range("E18").Value = "  +6.99%  "
$ws.Range("D19").Value = "3.018.52"
$ws.Range("E19").Value = "  +3.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "458.55"
$ws.Range("E20").Value = "  +6.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.96"
$ws.Range("E21").Value = "  +6.99%  "
$ws.Range("E22").Value = "  +5.99%  "
$ws.Range("E23").Value = "  +8.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.45"
$ws.Range("E24").Value = "  +4.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  +14.30%  "
$ws.Range("E26").Value = "  +3.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.67"
$ws.Range("E27").Value = "  +4.00%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.10"
$ws.Range("E29").Value = "  +16.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.38"
$ws.Range("E30").Value = "  +18.53%  "
$ws.Range("E31").Value = "  -6.40%  "
$ws.Range("E32").Value = "  +4.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.35"
$ws.Range("E33").Value = "  +6.50%  "
$ws.Range("E34").Value = "  +4.95%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.994"
$ws.Range("E36").Value = "  +4.06%  "
$ws.Range("E37").Value = "  +8.38%  "
$ws.Range("E38").Value = "  +15.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.03"
$ws.Range("E39").Value = "  +4.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.07"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("E41").Value = "  +16.56%  "
$ws.Range("E42").Value = "  +7.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.56"
$ws.Range("E43").Value = "  +5.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.47"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "398.95"
$ws.Range("E45").Value = "  +14.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0362"
$ws.Range("E46").Value = "  +7.74%  "
$ws.Range("D47").Value = "2.806.65"
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.71"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.01"
$ws.Range("E50").Value = "  +11.07%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.107"
$ws.Range("E51").Value = "  +4.48%  "
